$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Julio de 2020 a las 13:07"

# Row 14 - Iran
$ws.Range("B14").Value = 273788
$ws.Range("C14").Value = 2182
$ws.Range("D14").Value = 237788
$ws.Range("E14").Value = 21812
$ws.Range("G14").Value = 209
$ws.Range("H14").Value = 14188

# Row 40 - Emiratos Arabes Unidos
$ws.Range("B40").Value = 56922
$ws.Range("C40").Value = 211
$ws.Range("D40").Value = 49269
$ws.Range("E40").Value = 7314
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = 339

# Row 49 - Rumania
$ws.Range("B49").Value = 37458
$ws.Range("C49").Value = 767
$ws.Range("D49").Value = 22617
$ws.Range("E49").Value = 12815
$ws.Range("G49").Value = 17
$ws.Range("H49").Value = 2026

# Row 54 - Suiza
$ws.Range("B54").Value = 33591
$ws.Range("C54").Value = 99
$ws.Range("E54").Value = 1722

# Row 65 - Nepal
$ws.Range("B65").Value = 17658
$ws.Range("C65").Value = 156
$ws.Range("D65").Value = 11695
$ws.Range("E65").Value = 5923

# Row 113 - Sri Lanka
$ws.Range("B113").Value = 2711
$ws.Range("C113").Value = 7
$ws.Range("E113").Value = 665

# Row 125 - Islandia
$ws.Range("B125").Value = 1930
$ws.Range("C125").Value = 8
$ws.Range("D125").Value = 1907
$ws.Range("E125").Value = 13

# Row 162 - Vietnam
$ws.Range("B162").Value = 383
$ws.Range("C162").Value = 1
$ws.Range("E162").Value = 26

# Row 164 - Mauricio
$ws.Range("D164").Value = 332
$ws.Range("E164").Value = 1
